$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The top two rows (merged title "«Туризм и гостеприимство»..." and the
# "422" / "Заочное" group row) are removed entirely, shifting the header
# row and the two data rows up by two rows.
$ws.Rows("1:2").Delete()

# Re-anchor the AutoFilter on the new header row (was A3:G3, now A1:G1).
$ws.AutoFilterMode = $false
[void]$ws.Range("A1:G1").AutoFilter()

# Keep the workbook-level _FilterDatabase defined name in sync with the
# AutoFilter range.
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "='422'!`$A`$1:`$G`$1"
    }
}

# Restore the on-screen selection to match the post-edit layout.
[void]$ws.Range("A1:C2").Select()
